$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Hlavní město Praha" column (column B). This shifts every
# subsequent region column (old C:O) one place to the left (new B:N),
# matching the updated dataset/extrapolation values in the diff.
$ws.Range("B1").EntireColumn.Delete()
